# Update "想去人数" (interest counts) figures in the 展览 (Exhibition) sheet
# and mirror the same updates in the 全部类型 (All Types) sheet, which lists
# the same events one row further down.

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F6").Value = 9996
$wsExhibition.Range("F10").Value = 5525
$wsExhibition.Range("F17").Value = 296
$wsExhibition.Range("F18").Value = 588
$wsExhibition.Range("F22").Value = 1511

$wsAllTypes = $wb.Worksheets.Item("全部类型")
$wsAllTypes.Range("F7").Value = 9996
$wsAllTypes.Range("F11").Value = 5525
$wsAllTypes.Range("F18").Value = 296
$wsAllTypes.Range("F19").Value = 588
$wsAllTypes.Range("F23").Value = 1511
